$d = $word.ActiveDocument

# Map of 1-based paragraph index -> replacement WordprocessingML for the whole paragraph
# (drops w:proofErr spell/grammar markers and merges runs per the authoritative edit).
$fragments = @{}
$fragments[35] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>Fn er blevet beregnet og præsenteret for brugeren</w:t></w:r></w:p>'
$fragments[33] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>Fdim bliver beregnet og præsenteret for brugeren</w:t></w:r></w:p>'
$fragments[31] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>valgtProfil</w:t></w:r></w:p>'
$fragments[29] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>beregnFn(Fdim,vinkel)</w:t></w:r></w:p>'
$fragments[27] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Titel"/><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>PTE  - OC03: beregnFn</w:t></w:r></w:p>'
$fragments[22] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>Fdim er blevet bereg</w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>net og præsenteret til brugeren.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$fragments[21] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>dimensionkraft.getFdim er blevet kaldt</w:t></w:r></w:p>'
$fragments[20] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>Dimensionkraft dimensionkraft er blevet skabt med vægt som parameter.</w:t></w:r></w:p>'
$fragments[15] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>valgtProfil</w:t></w:r></w:p>'
$fragments[13] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>getFdim</w:t></w:r></w:p>'
$fragments[11] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Titel"/><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>PTE  - OC02: beregnFdim</w:t></w:r></w:p>'
$fragments[3] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>vælgProfil</w:t></w:r></w:p>'
$fragments[1] = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Titel"/><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>PTE</w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t xml:space="preserve">  - OC0</w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>vælgProfil</w:t></w:r></w:p>'

# Apply highest paragraph index first so earlier replacements do not disturb
# the character offsets backing not-yet-processed paragraphs.
$indices = $fragments.Keys | Sort-Object -Descending
foreach ($i in $indices) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    $rng.InsertXML($fragments[$i])
}

Write-Output "done"
